$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.18227207660675
$ws.Range("B1").Value = 1.662708640098572
$ws.Range("C1").Value = 2.920659303665161
$ws.Range("D1").Value = 1.508668541908264
$ws.Range("E1").Value = 0.8206518888473511
